$wb = $excel.ActiveWorkbook

# --- Project_Title sheet: engagement string update ---
$wsProj = $wb.Worksheets.Item("Project_Title")
$wsProj.Activate()

$oldEngagement = $wsProj.Range("A2").Value()
$wsProj.Range("D2").Value = $oldEngagement
$wsProj.Range("A2").Value = "GE Healthcare-GE Healthcare Bio-Sciences AB-FVA-101397"

# Header cells become bold
$wsProj.Range("A1:B1").Font.Bold = $true

# Column widths (best-fit-like) for the now-wider content
$wsProj.Columns.Item(1).ColumnWidth = 48.721354166666664
$wsProj.Columns.Item(2).ColumnWidth = 14.944010416666666
$wsProj.Columns.Item(4).ColumnWidth = 33.276041666666664

# Page setup: portrait
$wsProj.PageSetup.Orientation = 1

# Selection moves to A2
$wsProj.Range("A2").Select() | Out-Null

# --- Update_Hours sheet: selection change, no longer the active tab ---
$wsHours = $wb.Worksheets.Item("Update_Hours")
$wsHours.Activate()
$wsHours.Range("J22").Select() | Out-Null

# --- Error_Message sheet: becomes the active/selected tab ---
$wsErr = $wb.Worksheets.Item("Error_Message")
$wsErr.Activate()
$wsErr.Range("I22").Select() | Out-Null
